$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.433.67'
$ws.Range("E2").Value = '  +1.86%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.842.38'
$ws.Range("E3").Value = '  +1.41%  '
$ws.Range("E4").Value = '  +1.43%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.72'
$ws.Range("E5").Value = '  +1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.014'
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4759'
$ws.Range("E7").Value = '  +1.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3710'
$ws.Range("E8").Value = '  +0.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07468'
$ws.Range("E9").Value = '  +1.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8876'
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.52'
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.858.62'
$ws.Range("E12").Value = '  +2.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07382'
$ws.Range("E13").Value = '  +4.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.486'
$ws.Range("E14").Value = '  +1.77%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.37'
$ws.Range("E15").Value = '  +1.61%  '
$ws.Range("E16").Value = '  +1.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008856'
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("E19").Value = '  +1.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.85'
$ws.Range("E20").Value = '  +0.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.452.24'
$ws.Range("E21").Value = '  +1.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.352'
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("E23").Value = '  +1.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.074.91'
$ws.Range("E24").Value = '  +1.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.906'
$ws.Range("E25").Value = '  +0.78%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.58'
$ws.Range("E26").Value = '  +1.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.65'
$ws.Range("E27").Value = '  +1.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.167'
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.291'
$ws.Range("E29").Value = '  -0.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.22'
$ws.Range("E30").Value = '  +1.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08983'
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7611'
$ws.Range("E32").Value = '  -1.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.181'
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.570'
$ws.Range("E34").Value = '  +1.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.952'
$ws.Range("E35").Value = '  +1.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.014'
$ws.Range("E36").Value = '  +1.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.106'
$ws.Range("E37").Value = '  +1.88%  '
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01971'
$ws.Range("E39").Value = '  +0.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.002'
$ws.Range("E40").Value = '  +1.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.324'
$ws.Range("E41").Value = '  +0.58%  '
$ws.Range("E42").Value = '  +0.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.379'
$ws.Range("E43").Value = '  +1.78%  '
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.564'
$ws.Range("E45").Value = '  +1.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4982'
$ws.Range("E46").Value = '  +1.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.59'
$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.37'
$ws.Range("E49").Value = '  +1.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.686'
$ws.Range("E50").Value = '  +0.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06324'
